# The deck ships two theme parts:
#   ppt/theme/theme2.xml ("Integral" colours) -- linked from the slide master,
#       i.e. this is the theme that is actually visible on the slides.
#   ppt/theme/theme1.xml ("Office Theme" colours) -- linked from the notes
#       master only.
#
# The target edit swaps the two themes' content in full: theme1.xml ends up
# holding the "Integral" colour scheme and theme2.xml ends up holding the
# stock "Office Theme" colour scheme (the font scheme / format scheme
# portions of both theme parts are identical already, so only the 12
# clrScheme entries actually change).
#
# PowerPoint's automation model exposes a slide master's live theme palette
# through SlideMaster.Theme.ThemeColorScheme.Colors(1-12).RGB, in the fixed
# order dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -- so we recolour the
# visible theme (theme2.xml) to the stock Office palette that used to live
# in theme1.xml.

function ToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Stock "Office Theme" colour scheme, in Colors(1..12) order.
$officeThemeHex = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $tcs.Colors($i).RGB = ToOleRgb $officeThemeHex[$i - 1]
}
